# Insert a brand-new, most-recent log entry at the top of the September
# "Others" transaction list on the "2024" sheet.
#
# The list lives in columns R (September_Details) / S (September_Date) and
# runs from row 47 (most recent) down through row 181 (oldest), followed by
# blank padding rows up to the next category header ("Broadband", which
# currently sits on row 190).
#
# Inserting a new row above row 47 pushes that whole list (and the
# "Broadband" header below it) down by one row, and a fresh entry is
# written into the new row 47.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$ws.Rows("47").Insert()

$ws.Range("R47").Value = "credit icici"
$ws.Range("S47").Value = "2024-09-22 15:31:31"
